# Updated cryptos list on Thu Feb  8 13:26:15 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.726.32"
$ws.Range("E2").Value = "  +3.81%  "
$ws.Range("D3").Value = "2.419.05"
$ws.Range("E3").Value = "  +2.00%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "'317.23"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.66%  "
$ws.Range("D6").Value = "'101.66"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.57%  "
$ws.Range("D7").Value = "'0.512"
$ws.Range("D7").Style = "Normal"
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").Value = "'0.527"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +9.79%  "
$ws.Range("D10").Value = "'35.32"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.99%  "
$ws.Range("D11").Value = "'0.0799"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.58%  "
$ws.Range("E12").Value = "  -1.92%  "
$ws.Range("D13").Value = "'18.49"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.83%  "
$ws.Range("D14").Value = "'6.88"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.06%  "
$ws.Range("D15").Value = "2.797.53"
$ws.Range("E15").Value = "  +2.17%  "
$ws.Range("D16").Value = "2.443.68"
$ws.Range("E16").Value = "  +3.07%  "
$ws.Range("D17").Value = "'0.826"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.61%  "
$ws.Range("D18").Value = "44.535.51"
$ws.Range("E18").Value = "  +3.33%  "
$ws.Range("D19").Value = "'12.23"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.20%  "
$ws.Range("E20").Value = "  +1.17%  "
$ws.Range("D21").Value = "0.0₃0918"
$ws.Range("D22").Value = "'68.54"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.86%  "
$ws.Range("D23").Value = "'242.73"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.01%  "
$ws.Range("D24").Value = "'2.26"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.26%  "
$ws.Range("E25").Value = "  +2.29%  "
$ws.Range("E26").Value = "  -0.04%  "
$ws.Range("D27").Value = "'25.14"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.87%  "
$ws.Range("E28").Value = "  -3.46%  "
$ws.Range("D29").Value = "'9.49"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.45%  "
$ws.Range("D30").Value = "'33.37"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.52%  "
$ws.Range("D31").Value = "'48.27"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.74%  "
$ws.Range("D32").Value = "'0.125"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +13.70%  "
$ws.Range("D33").Value = "'19.51"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +10.75%  "
$ws.Range("D34").Value = "'5.17"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.01%  "
$ws.Range("E35").Value = "  +0.20%  "
$ws.Range("D36").Value = "'0.0762"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.49%  "
$ws.Range("E37").Value = "  +2.43%  "
$ws.Range("D38").Value = "'4.45"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.42%  "
$ws.Range("D39").Value = "'126.03"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.97%  "
$ws.Range("E40").Value = "  -0.27%  "
$ws.Range("E41").Value = "  +1.76%  "
$ws.Range("E42").Value = "  -3.84%  "
$ws.Range("E43").Value = "  +1.07%  "
$ws.Range("D44").Value = "'0.0288"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.47%  "
$ws.Range("D45").Value = "1.933.96"
$ws.Range("E45").Value = "  +0.22%  "
$ws.Range("E46").Value = "  -0.60%  "
$ws.Range("D47").Value = "'2.92"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +6.79%  "
$ws.Range("D48").Value = "'9.12"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.77%  "
$ws.Range("D49").Value = "'1.76"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +16.93%  "
$ws.Range("D50").Value = "'75.76"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.08%  "
$ws.Range("D51").Value = "'53.55"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.67%  "
